$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet tracks weekly price quotes for "Sandia" (watermelon), newest
# entries on top. A new week's data was added on top (rows 83:84), pushing
# the previous two rows (old 83:84, quality "Extra"/"Primera" @ 45233) down
# to become rows 85:86, and every row after that shifts down by two as well.
#
# Duplicate rows 83:84 into newly inserted rows 85:86 (this also shifts the
# old rows 85:95 down to 87:97), then overwrite rows 83:84 with the new
# week's figures.
$ws.Range("A83:R84").Copy()
$ws.Range("A85:R86").Insert()

# New week: row 83 ("Extra") - same volume/date as before except the date
# and price columns move.
$ws.Cells.Item(83, 4).Value = 45258
$ws.Cells.Item(83, 11).Value = 540
$ws.Cells.Item(83, 12).Value = 550
$ws.Cells.Item(83, 13).Value = 542
$ws.Cells.Item(83, 16).Value = 542

# New week: row 84 ("Primera")
$ws.Cells.Item(84, 4).Value = 45258
$ws.Cells.Item(84, 10).Value = 550
$ws.Cells.Item(84, 11).Value = 540
$ws.Cells.Item(84, 12).Value = 550
$ws.Cells.Item(84, 13).Value = 545
$ws.Cells.Item(84, 16).Value = 545
